$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Rubel Hossain*" row (row 23); all rows below shift up by one.
$ws.Rows.Item(23).Delete()

# The "Total Players" summary row (now row 35) needs its count updated from 34 to 33.
$ws.Range("A35").Value = "Total Players = 33"

# Update the query-table defined name range to reflect the new row extent (A1:P36 -> A1:P35).
$n = $wb.Names.Item(1)
$n.RefersTo = "=Sheet1!`$A`$1:`$P`$35"
